# Update cryptos list with latest price/volume data
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.290.01'
$ws.Range("E2").Value = '  +2.41%  '
$ws.Range("D3").Value = '3.170.43'
$ws.Range("E3").Value = '  +3.95%  '
$ws.Range("E4").Value = '  +0.01%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '578.25'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +3.91%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '150.70'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +6.32%  '
$ws.Range("E7").Value = '  +0.07%  '
$ws.Range("D8").Value = '3.169.59'
$ws.Range("E8").Value = '  +3.96%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.530'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +2.00%  '
$ws.Range("E10").Value = '  +6.33%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '6.19'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +0.69%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.502'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +4.83%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.0000269'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +16.59%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '37.42'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +6.60%  '
$ws.Range("D15").Value = '3.693.13'
$ws.Range("E15").Value = '  +3.96%  '
$ws.Range("D16").Value = '65.366.79'
$ws.Range("E16").Value = '  +2.45%  '
$ws.Range("B17").Value = 'WrappedEther'
$ws.Range("C17").Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range("D17").Value = '3.169.76'
$ws.Range("E17").Value = '  +4.01%  '
$ws.Range("B18").Value = 'Polkadot'
$ws.Range("C18").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '7.17'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +6.09%  '
$ws.Range("E19").Value = '  +1.41%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '512.08'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +5.33%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '14.88'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +5.07%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '0.724'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.99%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '15.31'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +5.13%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '7.81'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +4.03%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '85.12'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +3.20%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '9.07'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +11.94%  '
$ws.Range("E28").Value = '  +5.31%  '
$ws.Range("E29").Value = '  +7.67%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.83'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +15.29%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '27.87'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +5.63%  '
$ws.Range("E32").Value = '  -0.05%  '
$ws.Range("E33").Value = '  +4.24%  '
$ws.Range("E34").Value = '  +11.38%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '6.58'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +6.27%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '55.29'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  -0.13%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.0911'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  +12.00%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '474.53'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +7.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '0.0425'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +3.90%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '3.06'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +10.90%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '8.68'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +4.28%  '
$ws.Range("D42").Value = '3.065.60'
$ws.Range("E42").Value = '  +0.68%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '0.118'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +2.24%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '2.47'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +10.84%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.285'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +4.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '29.10'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +5.29%  '
$ws.Range("E47").Value = '  +17.61%  '
$ws.Range("E49").Value = '  +1.25%  '
$ws.Range("E50").Value = '  +8.07%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '121.57'
$ws.Range("D51").Style = "Normal"
